# This script updates the LR-pairs (Ntn4-Unc5a) worksheet with refreshed TPM-based
# NATMI values: a new "Resolving-Mac" cluster is introduced (adding 4 new data rows,
# rows 14-17), and the numeric statistics for the existing Sending/Target cluster
# combinations are recomputed against the new cluster set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Ntn4-Unc5a)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntn4"
$ws.Cells.Item(2, 3).Value = "Unc5a"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.136976666666667
$ws.Cells.Item(2, 8).Value = 3.41093
$ws.Cells.Item(2, 9).Value = 0.03386532673582325
$ws.Cells.Item(2, 10).Value = 0.03386532673582325
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.294217333333333
$ws.Cells.Item(2, 14).Value = 3.882652
$ws.Cells.Item(2, 15).Value = 0.1864098899142058
$ws.Cells.Item(2, 16).Value = 0.1864098899142058
$ws.Cells.Item(2, 17).Value = 1.471494909595556
$ws.Cells.Item(2, 18).Value = 13.24345418636
$ws.Cells.Item(2, 19).Value = 0.006312831828733424
$ws.Cells.Item(2, 20).Value = 0.006312831828733424

# Row 3: ECs -> FAPs (Ntn4-Unc5a)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntn4"
$ws.Cells.Item(3, 3).Value = "Unc5a"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.136976666666667
$ws.Cells.Item(3, 8).Value = 3.41093
$ws.Cells.Item(3, 9).Value = 0.03386532673582325
$ws.Cells.Item(3, 10).Value = 0.03386532673582325
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.626459
$ws.Cells.Item(3, 14).Value = 4.879377
$ws.Cells.Item(3, 15).Value = 0.2342636243010983
$ws.Cells.Item(3, 16).Value = 0.2342636243010983
$ws.Cells.Item(3, 17).Value = 1.84924593229
$ws.Cells.Item(3, 18).Value = 16.64321339061
$ws.Cells.Item(3, 19).Value = 0.007933414179274838
$ws.Cells.Item(3, 20).Value = 0.007933414179274836

# Row 4: ECs -> MuSCs (Ntn4-Unc5a)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntn4"
$ws.Cells.Item(4, 3).Value = "Unc5a"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.136976666666667
$ws.Cells.Item(4, 8).Value = 3.41093
$ws.Cells.Item(4, 9).Value = 0.03386532673582325
$ws.Cells.Item(4, 10).Value = 0.03386532673582325
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.488917666666667
$ws.Cells.Item(4, 14).Value = 10.466753
$ws.Cells.Item(4, 15).Value = 0.5025189675740148
$ws.Cells.Item(4, 16).Value = 0.5025189675740148
$ws.Cells.Item(4, 17).Value = 3.966817978921112
$ws.Cells.Item(4, 18).Value = 35.70136181029
$ws.Cells.Item(4, 19).Value = 0.01701796902784258
$ws.Cells.Item(4, 20).Value = 0.01701796902784258

# Row 5: ECs -> Resolving-Mac (Ntn4-Unc5a)
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Ntn4"
$ws.Cells.Item(5, 3).Value = "Unc5a"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.136976666666667
$ws.Cells.Item(5, 8).Value = 3.41093
$ws.Cells.Item(5, 9).Value = 0.03386532673582325
$ws.Cells.Item(5, 10).Value = 0.03386532673582325
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.5332636666666667
$ws.Cells.Item(5, 14).Value = 1.599791
$ws.Cells.Item(5, 15).Value = 0.07680751821068107
$ws.Cells.Item(5, 16).Value = 0.07680751821068106
$ws.Cells.Item(5, 17).Value = 0.6063083461811112
$ws.Cells.Item(5, 18).Value = 5.45677511563
$ws.Cells.Item(5, 19).Value = 0.002601111699972409
$ws.Cells.Item(5, 20).Value = 0.002601111699972409

# Row 6: FAPs -> ECs (Ntn4-Unc5a)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ntn4"
$ws.Cells.Item(6, 3).Value = "Unc5a"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 12.206517
$ws.Cells.Item(6, 8).Value = 36.619551
$ws.Cells.Item(6, 9).Value = 0.3635762268748239
$ws.Cells.Item(6, 10).Value = 0.3635762268748239
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.294217333333333
$ws.Cells.Item(6, 14).Value = 3.882652
$ws.Cells.Item(6, 15).Value = 0.1864098899142058
$ws.Cells.Item(6, 16).Value = 0.1864098899142058
$ws.Cells.Item(6, 17).Value = 15.797885881028
$ws.Cells.Item(6, 18).Value = 142.180972929252
$ws.Cells.Item(6, 19).Value = 0.06777420442715824
$ws.Cells.Item(6, 20).Value = 0.06777420442715824

# Row 7: FAPs -> FAPs (Ntn4-Unc5a)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ntn4"
$ws.Cells.Item(7, 3).Value = "Unc5a"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 12.206517
$ws.Cells.Item(7, 8).Value = 36.619551
$ws.Cells.Item(7, 9).Value = 0.3635762268748239
$ws.Cells.Item(7, 10).Value = 0.3635762268748239
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.626459
$ws.Cells.Item(7, 14).Value = 4.879377
$ws.Cells.Item(7, 15).Value = 0.2342636243010983
$ws.Cells.Item(7, 16).Value = 0.2342636243010983
$ws.Cells.Item(7, 17).Value = 19.853399433303
$ws.Cells.Item(7, 18).Value = 178.680594899727
$ws.Cells.Item(7, 19).Value = 0.08517268461741462
$ws.Cells.Item(7, 20).Value = 0.0851726846174146

# Row 8: FAPs -> MuSCs (Ntn4-Unc5a)
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Ntn4"
$ws.Cells.Item(8, 3).Value = "Unc5a"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 12.206517
$ws.Cells.Item(8, 8).Value = 36.619551
$ws.Cells.Item(8, 9).Value = 0.3635762268748239
$ws.Cells.Item(8, 10).Value = 0.3635762268748239
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.488917666666667
$ws.Cells.Item(8, 14).Value = 10.466753
$ws.Cells.Item(8, 15).Value = 0.5025189675740148
$ws.Cells.Item(8, 16).Value = 0.5025189675740148
$ws.Cells.Item(8, 17).Value = 42.587532809767
$ws.Cells.Item(8, 18).Value = 383.287795287903
$ws.Cells.Item(8, 19).Value = 0.1827039501635923
$ws.Cells.Item(8, 20).Value = 0.1827039501635923

# Row 9: FAPs -> Resolving-Mac (Ntn4-Unc5a)
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Ntn4"
$ws.Cells.Item(9, 3).Value = "Unc5a"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 12.206517
$ws.Cells.Item(9, 8).Value = 36.619551
$ws.Cells.Item(9, 9).Value = 0.3635762268748239
$ws.Cells.Item(9, 10).Value = 0.3635762268748239
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.5332636666666667
$ws.Cells.Item(9, 14).Value = 1.599791
$ws.Cells.Item(9, 15).Value = 0.07680751821068107
$ws.Cells.Item(9, 16).Value = 0.07680751821068106
$ws.Cells.Item(9, 17).Value = 6.509292012649
$ws.Cells.Item(9, 18).Value = 58.583628113841
$ws.Cells.Item(9, 19).Value = 0.02792538766665875
$ws.Cells.Item(9, 20).Value = 0.02792538766665875

# Row 10: MuSCs -> ECs (Ntn4-Unc5a)
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Ntn4"
$ws.Cells.Item(10, 3).Value = "Unc5a"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 20.10609633333333
$ws.Cells.Item(10, 8).Value = 60.318289
$ws.Cells.Item(10, 9).Value = 0.5988685095064435
$ws.Cells.Item(10, 10).Value = 0.5988685095064435
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.294217333333333
$ws.Cells.Item(10, 14).Value = 3.882652
$ws.Cells.Item(10, 15).Value = 0.1864098899142058
$ws.Cells.Item(10, 16).Value = 0.1864098899142058
$ws.Cells.Item(10, 17).Value = 26.02165838026978
$ws.Cells.Item(10, 18).Value = 234.194925422428
$ws.Cells.Item(10, 19).Value = 0.1116350129301807
$ws.Cells.Item(10, 20).Value = 0.1116350129301807

# Row 11: MuSCs -> FAPs (Ntn4-Unc5a)
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Ntn4"
$ws.Cells.Item(11, 3).Value = "Unc5a"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 20.10609633333333
$ws.Cells.Item(11, 8).Value = 60.318289
$ws.Cells.Item(11, 9).Value = 0.5988685095064435
$ws.Cells.Item(11, 10).Value = 0.5988685095064435
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.626459
$ws.Cells.Item(11, 14).Value = 4.879377
$ws.Cells.Item(11, 15).Value = 0.2342636243010983
$ws.Cells.Item(11, 16).Value = 0.2342636243010983
$ws.Cells.Item(11, 17).Value = 32.701741336217
$ws.Cells.Item(11, 18).Value = 294.315672025953
$ws.Cells.Item(11, 19).Value = 0.1402931075167762
$ws.Cells.Item(11, 20).Value = 0.1402931075167762

# Row 12: MuSCs -> MuSCs (Ntn4-Unc5a)
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Ntn4"
$ws.Cells.Item(12, 3).Value = "Unc5a"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 20.10609633333333
$ws.Cells.Item(12, 8).Value = 60.318289
$ws.Cells.Item(12, 9).Value = 0.5988685095064435
$ws.Cells.Item(12, 10).Value = 0.5988685095064435
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.488917666666667
$ws.Cells.Item(12, 14).Value = 10.466753
$ws.Cells.Item(12, 15).Value = 0.5025189675740148
$ws.Cells.Item(12, 16).Value = 0.5025189675740148
$ws.Cells.Item(12, 17).Value = 70.14851470506856
$ws.Cells.Item(12, 18).Value = 631.3366323456171
$ws.Cells.Item(12, 19).Value = 0.3009427851097671
$ws.Cells.Item(12, 20).Value = 0.3009427851097671

# Row 13: MuSCs -> Resolving-Mac (Ntn4-Unc5a)
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Ntn4"
$ws.Cells.Item(13, 3).Value = "Unc5a"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 20.10609633333333
$ws.Cells.Item(13, 8).Value = 60.318289
$ws.Cells.Item(13, 9).Value = 0.5988685095064435
$ws.Cells.Item(13, 10).Value = 0.5988685095064435
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.5332636666666667
$ws.Cells.Item(13, 14).Value = 1.599791
$ws.Cells.Item(13, 15).Value = 0.07680751821068107
$ws.Cells.Item(13, 16).Value = 0.07680751821068106
$ws.Cells.Item(13, 17).Value = 10.72185065306656
$ws.Cells.Item(13, 18).Value = 96.496655877599
$ws.Cells.Item(13, 19).Value = 0.04599760394971959
$ws.Cells.Item(13, 20).Value = 0.04599760394971959

# Row 14: Resolving-Mac -> ECs (Ntn4-Unc5a)
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Ntn4"
$ws.Cells.Item(14, 3).Value = "Unc5a"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.123884
$ws.Cells.Item(14, 8).Value = 0.371652
$ws.Cells.Item(14, 9).Value = 0.003689936882909407
$ws.Cells.Item(14, 10).Value = 0.003689936882909406
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.294217333333333
$ws.Cells.Item(14, 14).Value = 3.882652
$ws.Cells.Item(14, 15).Value = 0.1864098899142058
$ws.Cells.Item(14, 16).Value = 0.1864098899142058
$ws.Cells.Item(14, 17).Value = 0.1603328201226667
$ws.Cells.Item(14, 18).Value = 1.442995381104
$ws.Cells.Item(14, 19).Value = 0.0006878407281335103
$ws.Cells.Item(14, 20).Value = 0.0006878407281335102

# Row 15: Resolving-Mac -> FAPs (Ntn4-Unc5a)
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Ntn4"
$ws.Cells.Item(15, 3).Value = "Unc5a"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.123884
$ws.Cells.Item(15, 8).Value = 0.371652
$ws.Cells.Item(15, 9).Value = 0.003689936882909407
$ws.Cells.Item(15, 10).Value = 0.003689936882909406
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.626459
$ws.Cells.Item(15, 14).Value = 4.879377
$ws.Cells.Item(15, 15).Value = 0.2342636243010983
$ws.Cells.Item(15, 16).Value = 0.2342636243010983
$ws.Cells.Item(15, 17).Value = 0.201492246756
$ws.Cells.Item(15, 18).Value = 1.813430220804
$ws.Cells.Item(15, 19).Value = 0.0008644179876326549
$ws.Cells.Item(15, 20).Value = 0.0008644179876326547

# Row 16: Resolving-Mac -> MuSCs (Ntn4-Unc5a)
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Ntn4"
$ws.Cells.Item(16, 3).Value = "Unc5a"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.123884
$ws.Cells.Item(16, 8).Value = 0.371652
$ws.Cells.Item(16, 9).Value = 0.003689936882909407
$ws.Cells.Item(16, 10).Value = 0.003689936882909406
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 3.488917666666667
$ws.Cells.Item(16, 14).Value = 10.466753
$ws.Cells.Item(16, 15).Value = 0.5025189675740148
$ws.Cells.Item(16, 16).Value = 0.5025189675740148
$ws.Cells.Item(16, 17).Value = 0.4322210762173334
$ws.Cells.Item(16, 18).Value = 3.889989685956
$ws.Cells.Item(16, 19).Value = 0.001854263272812913
$ws.Cells.Item(16, 20).Value = 0.001854263272812913

# Row 17: Resolving-Mac -> Resolving-Mac (Ntn4-Unc5a)
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Ntn4"
$ws.Cells.Item(17, 3).Value = "Unc5a"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.123884
$ws.Cells.Item(17, 8).Value = 0.371652
$ws.Cells.Item(17, 9).Value = 0.003689936882909407
$ws.Cells.Item(17, 10).Value = 0.003689936882909406
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.5332636666666667
$ws.Cells.Item(17, 14).Value = 1.599791
$ws.Cells.Item(17, 15).Value = 0.07680751821068107
$ws.Cells.Item(17, 16).Value = 0.07680751821068106
$ws.Cells.Item(17, 17).Value = 0.06606283608133333
$ws.Cells.Item(17, 18).Value = 0.594565524732
$ws.Cells.Item(17, 19).Value = 0.000283414894330328
$ws.Cells.Item(17, 20).Value = 0.0002834148943303279
